# Apply weekly data refresh: update changed cells in rows 3-36
# Each entry: row, column index, new value
$changes = @(
  @(3,4,44364),
  @(3,10,100),
  @(3,11,19000),
  @(3,12,20000),
  @(3,13,19500),
  @(3,16,390),
  @(4,4,44364),
  @(4,10,100),
  @(4,11,19000),
  @(4,12,20000),
  @(4,13,19500),
  @(4,16,650),
  @(5,4,44385),
  @(5,8,'Española'),
  @(5,10,100),
  @(5,11,17000),
  @(5,12,18000),
  @(5,13,17500),
  @(5,14,'$/caja 30 unidades'),
  @(5,15,'Provincia de Limarí'),
  @(5,16,583),
  @(5,17,30),
  @(6,4,44484),
  @(6,10,220),
  @(6,11,8000),
  @(6,12,9000),
  @(6,13,8455),
  @(6,16,169),
  @(7,4,44484),
  @(7,10,220),
  @(7,11,7500),
  @(7,12,8000),
  @(7,13,7727),
  @(7,16,258),
  @(8,4,44397),
  @(8,8,'Española'),
  @(8,11,14000),
  @(8,12,15000),
  @(8,13,14500),
  @(8,14,'$/caja 30 unidades'),
  @(8,16,483),
  @(8,17,30),
  @(9,4,44383),
  @(9,8,'Argentina(o)'),
  @(9,10,50),
  @(9,11,17000),
  @(9,12,18000),
  @(9,13,17400),
  @(9,14,'$/caja 50 unidades'),
  @(9,16,348),
  @(9,17,50),
  @(10,4,44433),
  @(10,8,'Argentina(o)'),
  @(10,11,14000),
  @(10,12,15000),
  @(10,13,14500),
  @(10,14,'$/caja 50 unidades'),
  @(10,16,290),
  @(10,17,50),
  @(11,4,44335),
  @(11,8,'Española'),
  @(11,10,100),
  @(11,11,17000),
  @(11,12,18000),
  @(11,13,17500),
  @(11,14,'$/caja 30 unidades'),
  @(11,16,583),
  @(11,17,30),
  @(14,4,44497),
  @(14,8,'Argentina(o)'),
  @(14,10,180),
  @(14,11,6500),
  @(14,12,7000),
  @(14,13,6778),
  @(14,14,'$/caja 50 unidades'),
  @(14,16,136),
  @(14,17,50),
  @(15,4,44497),
  @(15,10,200),
  @(15,11,7000),
  @(15,12,7500),
  @(15,13,7250),
  @(15,16,242),
  @(16,4,44497),
  @(16,10,130),
  @(16,11,6000),
  @(16,12,6500),
  @(16,13,6192),
  @(16,16,155),
  @(17,4,44358),
  @(17,11,18000),
  @(17,12,20000),
  @(17,13,19000),
  @(17,16,380),
  @(18,8,'Española'),
  @(18,14,'$/caja 30 unidades'),
  @(18,16,633),
  @(18,17,30),
  @(19,4,44399),
  @(19,11,14000),
  @(19,12,15000),
  @(19,13,14500),
  @(19,16,483),
  @(20,4,44435),
  @(20,8,'Argentina(o)'),
  @(20,14,'$/caja 50 unidades'),
  @(20,16,290),
  @(20,17,50),
  @(21,4,44442),
  @(21,11,14500),
  @(21,12,15000),
  @(21,13,14750),
  @(21,16,492),
  @(22,4,44483),
  @(22,8,'Española'),
  @(22,10,450),
  @(22,11,11000),
  @(22,12,12000),
  @(22,13,11444),
  @(22,14,'$/caja 30 unidades'),
  @(22,16,381),
  @(22,17,30),
  @(23,4,44483),
  @(23,8,'Madrigal'),
  @(23,10,220),
  @(23,11,8000),
  @(23,12,8500),
  @(23,13,8273),
  @(23,14,'$/caja 40 unidades'),
  @(23,15,'Región de Coquimbo'),
  @(23,16,207),
  @(23,17,40),
  @(25,4,44420),
  @(25,11,14000),
  @(25,12,15000),
  @(25,13,14500),
  @(25,16,483),
  @(27,4,44463),
  @(27,8,'Argentina(o)'),
  @(27,11,9000),
  @(27,12,10000),
  @(27,13,9500),
  @(27,14,'$/caja 50 unidades'),
  @(27,16,190),
  @(27,17,50),
  @(28,4,44342),
  @(28,11,17000),
  @(28,12,18000),
  @(28,13,17500),
  @(28,16,583),
  @(29,4,44342),
  @(29,8,'Madrigal'),
  @(29,11,15000),
  @(29,12,16000),
  @(29,13,15500),
  @(29,14,'$/caja 40 unidades'),
  @(29,16,388),
  @(29,17,40),
  @(30,4,44376),
  @(30,8,'Española'),
  @(30,10,100),
  @(30,11,19000),
  @(30,12,20000),
  @(30,13,19500),
  @(30,14,'$/caja 30 unidades'),
  @(30,16,650),
  @(30,17,30),
  @(31,4,44421),
  @(31,10,100),
  @(31,11,14000),
  @(31,12,15000),
  @(31,13,14500),
  @(31,16,483),
  @(32,4,44426),
  @(32,8,'Madrigal'),
  @(32,10,50),
  @(32,11,12000),
  @(32,12,13000),
  @(32,13,12600),
  @(32,14,'$/caja 40 unidades'),
  @(32,16,315),
  @(32,17,40),
  @(33,4,44441),
  @(33,8,'Española'),
  @(33,11,13000),
  @(33,12,14000),
  @(33,13,13500),
  @(33,14,'$/caja 30 unidades'),
  @(33,16,450),
  @(33,17,30),
  @(34,4,44350),
  @(34,8,'Argentina(o)'),
  @(34,10,50),
  @(34,11,15000),
  @(34,12,16000),
  @(34,13,15600),
  @(34,14,'$/caja 50 unidades'),
  @(34,16,312),
  @(34,17,50),
  @(35,4,44350),
  @(35,8,'Española'),
  @(35,10,40),
  @(35,13,17500),
  @(35,14,'$/caja 30 unidades'),
  @(35,16,583),
  @(35,17,30),
  @(36,4,44447),
  @(36,8,'Española'),
  @(36,14,'$/caja 30 unidades'),
  @(36,16,483),
  @(36,17,30)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($chg in $changes) {
    $r = $chg[0]
    $c = $chg[1]
    $v = $chg[2]
    $ws.Cells.Item($r, $c).Value = $v
}

"Applied $($changes.Length) cell updates"
